# Weekly update: insert a new data row at the top of the Albahaca data block
# (row 171), pushing all existing data rows down by one. The previously
# last row (230) becomes row 231 unchanged; the new row 171 holds the
# newest weekly observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 171; this shifts rows 171:230 down to 172:231
# and Excel copies the formatting (incl. the date number format on column D)
# from the row being pushed down.
$ws.Rows.Item(171).Insert()

# Populate the newly inserted row 171 with the new weekly record.
$ws.Range("A171").Value = 10
$ws.Range("B171").Value = "Vega Modelo de Temuco"
$ws.Range("C171").Value = "La Araucanía"
$ws.Range("D171").Value = 44704
$ws.Range("E171").Value = 9
$ws.Range("F171").Value = 100112052
$ws.Range("G171").Value = "Albahaca"
$ws.Range("H171").Value = "Sin especificar"
$ws.Range("I171").Value = "Primera"
$ws.Range("J171").Value = 40
$ws.Range("K171").Value = 5000
$ws.Range("L171").Value = 5000
$ws.Range("M171").Value = 5000
$ws.Range("N171").Value = "$/paquete"
$ws.Range("O171").Value = "Región de Arica y Parinacota"
$ws.Range("P171").Value = 5000
$ws.Range("Q171").Value = 1
$ws.Range("R171").Value = "Hortaliza"
